$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 372, pushing the existing rows 372-400
# down to 374-402 (and the row styles/format along with them).
$ws.Rows("372:373").Insert()

# New row 372 (week of 2022-07-04, docena de paquetes / Provincia de Cautin)
$ws.Range("A372").Value = 10
$ws.Range("B372").Value = "Vega Modelo de Temuco"
$ws.Range("C372").Value = "La Araucanía"
$ws.Range("D372").Value = 44746
$ws.Range("E372").Value = 9
$ws.Range("F372").Value = 100114014
$ws.Range("G372").Value = "Betarraga"
$ws.Range("H372").Value = "Sin especificar"
$ws.Range("I372").Value = "Primera"
$ws.Range("J372").Value = 55
$ws.Range("K372").Value = 10000
$ws.Range("L372").Value = 10000
$ws.Range("M372").Value = 10000
$ws.Range("N372").Value = "`$/docena de paquetes"
$ws.Range("O372").Value = "Provincia de Cautín"
$ws.Range("P372").Value = 833
$ws.Range("Q372").Value = 12
$ws.Range("R372").Value = "Hortaliza"

# New row 373 (week of 2022-07-04, saco 25 kilos / Provincia de Cautin)
$ws.Range("A373").Value = 10
$ws.Range("B373").Value = "Vega Modelo de Temuco"
$ws.Range("C373").Value = "La Araucanía"
$ws.Range("D373").Value = 44746
$ws.Range("E373").Value = 9
$ws.Range("F373").Value = 100114014
$ws.Range("G373").Value = "Betarraga"
$ws.Range("H373").Value = "Sin especificar"
$ws.Range("I373").Value = "Primera"
$ws.Range("J373").Value = 110
$ws.Range("K373").Value = 8000
$ws.Range("L373").Value = 8000
$ws.Range("M373").Value = 8000
$ws.Range("N373").Value = "`$/saco 25 kilos"
$ws.Range("O373").Value = "Provincia de Cautín"
$ws.Range("P373").Value = 320
$ws.Range("Q373").Value = 25
$ws.Range("R373").Value = "Hortaliza"
